# Apply the "ISE" (Import Substitution Elasticities) acronym addition
# to the "Key to Variables" sheet, plus associated workbook/view state
# tweaks, per the commit "Add ISE data file for Import Substitution
# Elasticities (issue #183)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# --- Workbook-level calculation settings -------------------------------
# Turn on iterative calculation (calcPr iterate="1" iterateDelta="1e-5").
$excel.Iteration  = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.00001

# --- Insert the new "indst" row for ISE ---------------------------------
# A new row is inserted right above the existing "MHV" row (old row 148),
# shifting every subsequent row of the "indst" / "oes" blocks down by one.
$ws.Rows.Item(148).Insert()

$ws.Range("A148").Value = "indst"
$ws.Range("B148").Value = "ISE"
$ws.Range("C148").Value = "Import Substitution Elasticities"

# Copy the formatting (and value) from another "medium" importance cell
# (F4, style index 6 / fill "medium") onto the new F148 cell so the cell
# picks up the same named fill style used elsewhere for "medium". F4 sits
# well above the inserted row so it is unaffected by the row shift.
$ws.Range("F4").Copy($ws.Range("F148"))

$excel.CutCopyMode = $false

# --- View state: scroll / selection update ------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 131
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F148").Select()
